$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1911.8823
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1900.1333
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1900.1333
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2250.1333

# ALC row 52
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 5000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 15000
$ws.Range("M52").Value = ""
$ws.Range("N52").Value = -15320

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4547134.5
$ws.Range("I137").Value = 2274043
$ws.Range("J137").Value = 9093318
$ws.Range("K137").Value = 6822129
$ws.Range("L137").Value = 27279954
$ws.Range("M137").Value = -6819579
$ws.Range("N137").Value = -27285054

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1986.3125
$ws.Range("I141").Value = 1469.037
$ws.Range("J141").Value = 4779.6
$ws.Range("K141").Value = 4407.111
$ws.Range("L141").Value = 14338.8
$ws.Range("M141").Value = 772.8890000000001

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 373.33334
$ws.Range("I5").Value = 373.33334
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 373.33334
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -261.33334
$ws.Range("N5").Value = ""

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 826.5
$ws.Range("I97").Value = 739.8
$ws.Range("J97").Value = 913.2
$ws.Range("K97").Value = 739.8
$ws.Range("L97").Value = 913.2
$ws.Range("M97").Value = -243.8
$ws.Range("N97").Value = -1905.2

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 250008560
$ws.Range("I102").Value = 2105
$ws.Range("J102").Value = 500015000
$ws.Range("K102").Value = 2105
$ws.Range("L102").Value = 500015000
$ws.Range("M102").Value = -483

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1976.5
$ws.Range("I122").Value = 1976.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5929.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3479.5
$ws.Range("N122").Value = ""

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 373.33334
$ws.Range("I4").Value = 373.33334
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 373.33334
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -258.33334
$ws.Range("N4").Value = ""

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 336
$ws.Range("I22").Value = 336
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 336
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -163
$ws.Range("N22").Value = ""

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1225.5714
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 1394.75
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 1394.75
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -2296.75

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 40.8
$ws.Range("I7").Value = 36
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 36
$ws.Range("L7").Value = 60
$ws.Range("M7").Value = 77
$ws.Range("N7").Value = -286

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1796.931
$ws.Range("I31").Value = 1575.7959
$ws.Range("J31").Value = 3000.889
$ws.Range("K31").Value = 1575.7959
$ws.Range("L31").Value = 3000.889
$ws.Range("M31").Value = -1280.7959
$ws.Range("N31").Value = -3590.889

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1796.931
$ws.Range("I34").Value = 1575.7959
$ws.Range("J34").Value = 3000.889
$ws.Range("K34").Value = 1575.7959
$ws.Range("L34").Value = 3000.889
$ws.Range("M34").Value = -1373.7959
$ws.Range("N34").Value = -3404.889

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3161
$ws.Range("I105").Value = 1089.6666
$ws.Range("J105").Value = 4196.6665
$ws.Range("K105").Value = 1089.6666
$ws.Range("L105").Value = 4196.6665
$ws.Range("M105").Value = 657.3334
$ws.Range("N105").Value = -7690.6665

# CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 25314.182
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 25314.182
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 75942.546
$ws.Range("N9").Value = -76390.546

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 67700
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 101500
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 304500
$ws.Range("M17").Value = -131
$ws.Range("N17").Value = -304838

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1517.5385
$ws.Range("I34").Value = 540.6667
$ws.Range("J34").Value = 1810.6
$ws.Range("K34").Value = 1622.0001
$ws.Range("L34").Value = 5431.799999999999
$ws.Range("M34").Value = -1538.0001
$ws.Range("N34").Value = -5599.799999999999

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3386
$ws.Range("I46").Value = 633.3333
$ws.Range("J46").Value = 4136.727
$ws.Range("K46").Value = 1899.9999
$ws.Range("L46").Value = 12410.181
$ws.Range("M46").Value = -1808.9999
$ws.Range("N46").Value = -12592.181

# CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3570.9285
$ws.Range("I58").Value = 2499.5
$ws.Range("J58").Value = 3749.5
$ws.Range("K58").Value = 7498.5
$ws.Range("L58").Value = 11248.5
$ws.Range("M58").Value = -7370.5
$ws.Range("N58").Value = -11504.5

# CUL row 76
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 7633.3335
$ws.Range("I76").Value = 4000
$ws.Range("J76").Value = 7815
$ws.Range("K76").Value = 12000
$ws.Range("L76").Value = 23445
$ws.Range("M76").Value = -11617
$ws.Range("N76").Value = -24211

# CUL row 79
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 7633.3335
$ws.Range("I79").Value = 4000
$ws.Range("J79").Value = 7815
$ws.Range("K79").Value = 12000
$ws.Range("L79").Value = 23445
$ws.Range("M79").Value = -10674
$ws.Range("N79").Value = -26097

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 877.3333
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 895.375
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2686.125
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12766.125

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 3971.6667
$ws.Range("I133").Value = 3707.5
$ws.Range("J133").Value = 4500
$ws.Range("K133").Value = 11122.5
$ws.Range("L133").Value = 13500
$ws.Range("M133").Value = -6062.5
$ws.Range("N133").Value = -23620

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2342.1428
$ws.Range("I137").Value = 1629
$ws.Range("J137").Value = 4125
$ws.Range("K137").Value = 4887
$ws.Range("L137").Value = 12375
$ws.Range("M137").Value = 213
$ws.Range("N137").Value = -22575

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 28112.846
$ws.Range("I139").Value = 2324.1177
$ws.Range("J139").Value = 48040.5
$ws.Range("K139").Value = 6972.353099999999
$ws.Range("L139").Value = 144121.5
$ws.Range("M139").Value = -1832.353099999999
$ws.Range("N139").Value = -154401.5

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 9666
$ws.Range("I141").Value = 4123.75
$ws.Range("J141").Value = 16000
$ws.Range("K141").Value = 12371.25
$ws.Range("L141").Value = 48000
$ws.Range("M141").Value = -7191.25
$ws.Range("N141").Value = -58360

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2512.75
$ws.Range("I93").Value = 1620.4
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 1620.4
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -372.4000000000001
$ws.Range("N93").Value = -6496

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1865.7941
$ws.Range("I100").Value = 1837.7446
$ws.Range("J100").Value = 1928.5714
$ws.Range("K100").Value = 1837.7446
$ws.Range("L100").Value = 1928.5714
$ws.Range("M100").Value = -1296.7446
$ws.Range("N100").Value = -3010.5714

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4328.8887
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4422.857
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 13268.571
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -18168.571

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2526.9656
$ws.Range("I132").Value = 1965.9048
$ws.Range("J132").Value = 3999.75
$ws.Range("K132").Value = 5897.7144
$ws.Range("L132").Value = 11999.25
$ws.Range("M132").Value = -3367.7144
$ws.Range("N132").Value = -17059.25

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 46335
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 46335
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 46335
$ws.Range("N133").Value = -51395

# WVR row 80
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 31500
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 31500
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 31500
$ws.Range("N80").Value = -33496

# WVR row 83
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 31500
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 31500
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 94500
$ws.Range("N83").Value = -104484
